$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row values for the new dump entry
$ws.Range("B14").Value = "Dump20160301"
$ws.Range("C14").Value = "Added new field for REQUEST_ORDER_HEADER table"

# Move the active selection to C15, matching the updated sheet view
$ws.Range("C15").Select()
